$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New rows of data to append (rows 10-13), matching columns A..N
$rows = @(
    @(42613.759386574071, 8,  52, 47, 52, 40, 25980, 17407, 995, 181, 163, 13, 9, "Named"),
    @(42613.888055555559, 12, 54, 43, 54, 34, 17834, 17309, 946, 181, 146, 17, 9, "Named"),
    @(42614.886192129627, 32, 53, 44, 53, 4,  35605, 16366, 928, 168, 140, 20, 1, "Named"),
    @(42615.885671296295, 14, 55, 44, 55, 23, 19229, 16159, 837, 163, 130, 10, 3, "Named")
)

$startRow = 10
for ($i = 0; $i -lt $rows.Count; $i++) {
    $r = $startRow + $i
    $data = $rows[$i]
    for ($col = 0; $col -lt $data.Count; $col++) {
        $cell = $ws.Cells.Item($r, $col + 1)
        $cell.Value2 = $data[$col]
    }
    # Column A keeps the date number format used by the rest of the column
    $ws.Cells.Item($r, 1).NumberFormat = "m/d/yy h:mm"
}
